$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Relocate the old 4-row "header" block (old rows 12-15) down to rows 16-19 ---
# Row 15 (the last entry of that block) is removed first, while nothing exists below it yet,
# so the delete/shift only affects row 15 itself.
$ws.Rows(15).Delete()

# The remaining three old header rows (12,13,14) copy down to their new home (16,17,18),
# carrying both values and formatting.
$ws.Range("A12:O14").Copy($ws.Range("A16:O18"))

# Recreate what used to be row 15 at its new position, row 19: clone the formatting from the
# row above (now at 18) and re-enter its two CO-number / product values.
$ws.Range("A18:O18").Copy($ws.Range("A19:O19"))
$ws.Range("A19").Value = "3013696352"
$ws.Range("B19").Value = "TA5TVBPC"

# --- 2) Stamp the "detail" row formatting (as used on rows 8-10) onto the new detail rows ---
$ws.Range("A8:B8").Copy($ws.Range("A6:B6"))
$ws.Range("A8:B8").Copy($ws.Range("A7:B7"))
$ws.Range("A8:O8").Copy($ws.Range("A11:O11"))
$ws.Range("A8:O8").Copy($ws.Range("A12:O12"))
$ws.Range("A8:O8").Copy($ws.Range("A13:O13"))
$ws.Range("A8:O8").Copy($ws.Range("A14:O14"))

# --- 3) Stamp the "header" row formatting (as used on row 16, ex-row 12) onto the new rows 20-25 ---
$ws.Range("A16:B16").Copy($ws.Range("A20:B20"))
$ws.Range("A16:B16").Copy($ws.Range("A21:B21"))
$ws.Range("A16:B16").Copy($ws.Range("A22:B22"))
$ws.Range("A16:B16").Copy($ws.Range("A23:B23"))
$ws.Range("A16:B16").Copy($ws.Range("A24:B24"))
$ws.Range("A16:B16").Copy($ws.Range("A25:B25"))

# --- 4) Fill in all of the new CO-number / product values, top-to-bottom ---
$ws.Range("A5").Value = "3013696483"
$ws.Range("B5").Value = "TB7SX14CC"

$ws.Range("A6").Value = "3013696484"
$ws.Range("B6").Value = "TB7SX14CC"

$ws.Range("A7").Value = "3013696525"
$ws.Range("B7").Value = "TB7SX14CC"

$ws.Range("A8").Value = "3013696526"
$ws.Range("B8").Value = "TB7SX14CC"

$ws.Range("A9").Value = "3013696527"
$ws.Range("B9").Value = "TB7SX14CC"

$ws.Range("A10").Value = "3013696528"
$ws.Range("B10").Value = "TB7SX14CC"

$ws.Range("A11").Value = "3013696529"
$ws.Range("B11").Value = "TB7SX14CC"

$ws.Range("A12").Value = "3013696530"
$ws.Range("B12").Value = "TB7SX14CC"

$ws.Range("A13").Value = "3013696532"
$ws.Range("B13").Value = "TB7SX14CC"

$ws.Range("A14").Value = "3013696533"
$ws.Range("B14").Value = "TB7SX14CC"

$ws.Range("A20").Value = "3013696536"
$ws.Range("B20").Value = "TA5TVBPC"

$ws.Range("A21").Value = "3013696537"
$ws.Range("B21").Value = "TA5TVBPC"

$ws.Range("A22").Value = "3013696538"
$ws.Range("B22").Value = "TA5TVBPC"

$ws.Range("A23").Value = "3013696539"
$ws.Range("B23").Value = "TA5TVBPC"

$ws.Range("A24").Value = "3013696540"
$ws.Range("B24").Value = "TA5TVBPC"

$ws.Range("A25").Value = "3013696541"
$ws.Range("B25").Value = "TA5TVBPC"

# --- 5) Leave the selection on the last entered cell, like a user would after typing the list ---
$ws.Range("A25").Select()
